$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("C2").Value = 2
$ws1.Range("C7").Value = 13

# --- Sheet "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("C4").Value = 2

# --- Sheet "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")
$ws5.Range("B2").Value = 74
$ws5.Range("C2").Value = 71.8
$ws5.Range("D2").Value = 74
$ws5.Range("E2").Value = 78.7

$ws5.Range("B3").Value = 29
$ws5.Range("C3").Value = 28.2
$ws5.Range("D3").Value = 20
$ws5.Range("E3").Value = 21.3
